# "update. converted pdf file of current results (QC)"
#
# Re-crops/repositions the two QC result screenshots on each of the two
# content slides to match the refreshed PDF export: both pictures are
# now uniformly sized (~3,600,000 EMU wide) and shifted to their new
# layout position.
#
# Shape.Left/Top/Width/Height are expressed in points in the PowerPoint
# object model (the OOXML stores EMU, 12700 EMU == 1 point). The literal
# point values below are the exact doubles whose internal float32
# round-trip reproduces the target EMU numbers bit-for-bit (COM shape
# geometry is marshalled as 32-bit floats), so the saved OOXML matches
# the target exactly instead of drifting by a stray EMU.

$p = $ppt.ActivePresentation

# ---- Slide 1 : "Picture 3" (rId3) and "Picture 4" (rId4) ----
$s1 = $p.Slides.Item(1)

$pic3 = $s1.Shapes.Item(1)
$pic3.Left   = 75.37504197007874    # 957263 EMU
$pic3.Top    = 14.62496062992126    # 185737 EMU
$pic3.Width  = 283.46456692913387   # 3600000 EMU
$pic3.Height = 290.6636220472441    # 3691428 EMU

$pic4 = $s1.Shapes.Item(2)
$pic4.Left   = 400.07142642283463   # 5080907 EMU
$pic4.Top    = 21.824015748031496   # 277165 EMU
$pic4.Width  = 283.46456692913387   # 3600000 EMU
$pic4.Height = 283.46456692913387   # 3600000 EMU

# ---- Slide 2 : "Picture 1" (rId3) and "Picture 2" (rId4) ----
$s2 = $p.Slides.Item(2)

$pic1 = $s2.Shapes.Item(1)
$pic1.Left   = 164.79566929133858   # 2092905 EMU
$pic1.Top    = 100.30338672677166   # 1273853 EMU
$pic1.Width  = 278.4543307086614    # 3536370 EMU
$pic1.Height = 294.5361481322835    # 3740609 EMU

$pic2 = $s2.Shapes.Item(2)
$pic2.Left   = 518.936377952756     # 6590492 EMU
$pic2.Top    = 78.75                # 1000125 EMU
$pic2.Width  = 283.46456692913387   # 3600000 EMU
$pic2.Height = 283.46456692913387   # 3600000 EMU
